# add the NA's under duplicate_image_filename
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
